$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with two new columns P and Q (values 14 and 15),
# matching the style of the existing header cells (bold / bordered / centered).
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update the data columns for rows 2-25: swap I<->old K value, and M<->old O value.
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Add the two new data columns P and Q for rows 2-25, all with value 2.
$ws.Range("P2:Q25").Value = 2
